# Gantt Chart workbook touch-up:
#  - Fix the "Introduction to cozmo" start date (B10), which had been
#    mistakenly entered as a half-day decimal (43504.5 / 0.00 number format)
#    instead of a whole date like the rest of the Start Date column.
#  - Centre-align the Start Date and Days to Complete columns (B & C), which
#    also leaves those two columns narrower than before.
#  - Set up the page for printing (A4, portrait).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the bad date value/format in B10 and B21 --------------------------
$ws.Range("B10").Value = 43504
$ws.Range("B10").NumberFormat = "d-mmm"
$ws.Range("B21").NumberFormat = "d-mmm"

# --- Centre-align the Start Date and Days to Complete columns --------------
$ws.Range("B1:C19").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B21").HorizontalAlignment = -4108      # xlCenter

# --- Resize columns B and C (as a result of the formatting pass) -----------
$ws.Columns("B").ColumnWidth = 9.25
$ws.Columns("C").ColumnWidth = 16.95

# --- Selecting column C is what the author last did before saving ----------
$ws.Columns("C").Select()

# --- Configure the page for printing ----------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
